$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "I" column values for the parameters table (rows 19-24)
$ws.Range("I19").Value = 10.548500000000001
$ws.Range("I20").Value = 0.1769
$ws.Range("I21").Value = 0
$ws.Range("I23").Value = 0.0078
$ws.Range("I24").Value = 0.68200000000000005

# Update the selected cell to match the recorded view state
$ws.Range("G17").Select()
